$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 15.85673476965436
$ws.Range("C4").Value = 2844.333701422901
$ws.Range("D4").Value = 0.2342008361026222
$ws.Range("E4").Value = 0.3547441259492827
$ws.Range("F4").Value = 0.3237522793961501
$ws.Range("G4").Value = 0.4149427983271107
$ws.Range("H4").Value = 0.8208594541028789
$ws.Range("I4").Value = 0.6406939169195186
$ws.Range("J4").Value = 0.7048830455256553
$ws.Range("K4").Value = 0.4350921246976218
$ws.Range("L4").Value = 0.3185613570479362
$ws.Range("M4").Value = 0.3954179001576341
$ws.Range("N4").Value = 0.8442029262184203
$ws.Range("O4").Value = 0.8486165070739641
$ws.Range("P4").Value = 0.7306208688026333
$ws.Range("Q4").Value = 0.7175926563552628
$ws.Range("R4").Value = 0.3958050224149399
$ws.Range("S4").Value = 0.6103907608621331
$ws.Range("T4").Value = 0.3234441300793346
$ws.Range("U4").Value = 0.152638821955372
$ws.Range("V4").Value = 0.1222119178755813
$ws.Range("W4").Value = 0.147919988298509
$ws.Range("X4").Value = 0.04892618527123609
$ws.Range("Y4").Value = 0.1709991689434509
$ws.Range("Z4").Value = 0.0138612145644401
$ws.Range("AA4").Value = 0.03169347573938765
$ws.Range("AB4").Value = 0.1141949275098944
$ws.Range("AC4").Value = 0.1394472488074933
$ws.Range("AD4").Value = 0.2374767510971146
$ws.Range("AE4").Value = 0.2283336821612848
$ws.Range("AF4").Value = 0.4480568140853575
$ws.Range("AR4").Value = 15.68477249500009
$ws.Range("AS4").Value = 2867.047350529075
$ws.Range("AT4").Value = 3.103902235390978
$ws.Range("AU4").Value = 2.746860534078934
$ws.Range("AV4").Value = 3.170100932506231
$ws.Range("AW4").Value = 2.308876875234968
$ws.Range("AX4").Value = 0.9205233058978229
$ws.Range("AY4").Value = 0.9963167507817223
$ws.Range("AZ4").Value = 0.7290028985395032
$ws.Range("BA4").Value = 0.651559386944885
$ws.Range("BB4").Value = 3.676996858380262
$ws.Range("BC4").Value = 0.6613136221303231
$ws.Range("BD4").Value = 1.063889551893675
$ws.Range("BE4").Value = 0.9459959602198982
$ws.Range("BF4").Value = 1.292049486351956
$ws.Range("BG4").Value = 1.16831772991364
$ws.Range("BH4").Value = 1.66977828013312
$ws.Range("BI4").Value = 1.805114570961413
$ws.Range("BJ4").Value = 2.292281375720801
$ws.Range("BK4").Value = 2.416557100474932
$ws.Range("BL4").Value = 0.937192825752144
$ws.Range("BM4").Value = -0.01602713703362475
$ws.Range("BN4").Value = 0.06802224147861428
$ws.Range("BO4").Value = 0.04387100776631375
$ws.Range("BP4").Value = -0.1264984663865323
$ws.Range("BQ4").Value = 16.11778401744839
$ws.Range("BR4").Value = 2915.011355010498
